# Reorder columns in the consolidated hours report:
#   B (Total de horas) <-> C (Banco de horas)
#   D (Adicional noturno) <-> F (Horas S.T)
# Column E (Horas S.A) keeps its place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = $ws.Range("B5").Value()
$colC = $ws.Range("C5").Value()
$ws.Range("B5").Value = $colC
$ws.Range("C5").Value = $colB

$colD = $ws.Range("D5").Value()
$colF = $ws.Range("F5").Value()
$ws.Range("D5").Value = $colF
$ws.Range("F5").Value = $colD

# Columns D and E now share the same (narrower) width, matching column E's
# original width, since column D no longer holds the wider "Adicional
# noturno" header.
$ws.Range("D1:E1").ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Reset the sheet selection to the top header row instead of the data body.
$ws.Range("A1:F1").Select()
